{"js": "// Update the answers in the \"two-digit number divided by one-digit number\"\n// worksheet table. The table has 5 \"data\" rows (row indices 0, 4, 8, 12, 16)\n// each holding 5 division-fact cells; the remaining rows are blank spacer\n// rows. We address cells by (row, column) position rather than by searching\n// for the old text, because a couple of the new values coincide with old\n// values used elsewhere in the table (e.g. \"32\u00f77=4, 4\" is both a value being\n// replaced and a value being written), so naive global find/replace would\n// corrupt the result if applied sequentially.\n\nconst newValuesByRow = {\n  0: [\"23\u00f77=3, 2\", \"40\u00f78=5, 0\", \"97\u00f78=12, 1\", \"82\u00f78=10, 2\", \"98\u00f72=49, 0\"],\n  4: [\"91\u00f72=45, 1\", \"18\u00f72=9, 0\", \"69\u00f77=9, 6\", \"85\u00f77=12, 1\", \"67\u00f72=33, 1\"],\n  8: [\"37\u00f75=7, 2\", \"92\u00f77=13, 1\", \"32\u00f77=4, 4\", \"86\u00f78=10, 6\", \"17\u00f79=1, 8\"],\n  12: [\"83\u00f79=9, 2\", \"12\u00f74=3, 0\", \"43\u00f74=10, 3\", \"52\u00f73=17, 1\", \"94\u00f78=11, 6\"],\n  16: [\"84\u00f75=16, 4\", \"81\u00f77=11, 4\", \"13\u00f78=1, 5\", \"75\u00f75=15, 0\", \"50\u00f72=25, 0\"],\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const rowIndexStr of Object.keys(newValuesByRow)) {\n  const rowIndex = Number(rowIndexStr);\n  const newValues = newValuesByRow[rowIndex];\n\n  const cells = rows.items[rowIndex].cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  // Load each cell's first (only) paragraph so we can replace its range's\n  // text in place. Replacing via the paragraph's own range (instead of the\n  // whole cell body) keeps the existing run-level formatting (rPr: font,\n  // size) and paragraph formatting (pPr: alignment) untouched, so only the\n  // <w:t> content changes, matching the original diff.\n  const paragraphsByCol = [];\n  for (let c = 0; c < newValues.length; c++) {\n    const paragraphs = cells.items[c].body.paragraphs;\n    paragraphs.load(\"items\");\n    paragraphsByCol.push(paragraphs);\n  }\n  await context.sync();\n\n  for (let c = 0; c < newValues.length; c++) {\n    const firstParagraph = paragraphsByCol[c].items[0];\n    const range = firstParagraph.getRange();\n    range.insertText(newValues[c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the answers in the \"two-digit number divided by one-digit number\"\n# worksheet table. The table has 5 \"data\" rows (1-based COM row indices\n# 1, 5, 9, 13, 17) each holding 5 division-fact cells (columns 1-5); the\n# remaining rows are blank spacer rows. We address cells by (row, column)\n# position rather than by Find/Replace on the old text, because a couple of\n# the new values coincide with old values used elsewhere in the table (e.g.\n# \"32\u00f77=4, 4\" is both a value being replaced and a value being written), so\n# a naive global find-and-replace-all would corrupt the result if the rules\n# were applied sequentially.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValuesByRow = @{\n    1  = @(\"23\u00f77=3, 2\", \"40\u00f78=5, 0\", \"97\u00f78=12, 1\", \"82\u00f78=10, 2\", \"98\u00f72=49, 0\")\n    5  = @(\"91\u00f72=45, 1\", \"18\u00f72=9, 0\", \"69\u00f77=9, 6\", \"85\u00f77=12, 1\", \"67\u00f72=33, 1\")\n    9  = @(\"37\u00f75=7, 2\", \"92\u00f77=13, 1\", \"32\u00f77=4, 4\", \"86\u00f78=10, 6\", \"17\u00f79=1, 8\")\n    13 = @(\"83\u00f79=9, 2\", \"12\u00f74=3, 0\", \"43\u00f74=10, 3\", \"52\u00f73=17, 1\", \"94\u00f78=11, 6\")\n    17 = @(\"84\u00f75=16, 4\", \"81\u00f77=11, 4\", \"13\u00f78=1, 5\", \"75\u00f75=15, 0\", \"50\u00f72=25, 0\")\n}\n\nforeach ($rowIndex in $newValuesByRow.Keys) {\n    $newValues = $newValuesByRow[$rowIndex]\n    for ($col = 1; $col -le $newValues.Length; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        # Assigning to Range.Text replaces just the cell's text content (Word\n        # automatically keeps the end-of-cell marker out of the replace),\n        # preserving the existing run/paragraph formatting (font, size,\n        # alignment) so only the <w:t> content changes.\n        $cell.Range.Text = $newValues[$col - 1]\n    }\n}\n"}
